# 检查清单.xlsx - mark more chapters as studied/mastered and log dates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Select()

# Rows 12-17 (chapters 10-15): already had dates, now mark "learned" / "mastered" = 是
foreach ($r in 12..17) {
    $ws.Range("D$r").Value = "是"
    $ws.Range("E$r").Value = "是"
}

# Rows 18-20 (chapters 16-18): mark learned/mastered + add study date
$ws.Range("D18").Value = "是"
$ws.Range("E18").Value = "是"
$ws.Range("F18").Value = 43519

$ws.Range("D19").Value = "是"
$ws.Range("E19").Value = "是"
$ws.Range("F19").Value = 43519

$ws.Range("D20").Value = "是"
$ws.Range("E20").Value = "是"
$ws.Range("F20").Value = 43521

# Rows 21-24 (chapters 19-22): only the date gets recorded so far
$ws.Range("F21").Value = 43521
$ws.Range("F22").Value = 43521
$ws.Range("F23").Value = 43521
$ws.Range("F24").Value = 43521

# Row 28 (chapter 23, Java 继承): date recorded
$ws.Range("F28").Value = 43522

# Move the view to where today's work is happening and select the active cell
$ws.Range("F29").Select() | Out-Null
